# Add a new "2020" data column (Q) to the suicide-mortality-rate sheet,
# mirroring the formatting of the existing "2019" column (P), and select
# column T (matches the sheetView selection recorded in the target file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy the formatting of column P (2019) into column Q (2020) for the
#    header/data rows (4-14) so the new column matches the existing table
#    styling exactly.
$ws.Range("P4:P14").Copy() | Out-Null
$ws.Range("Q4:Q14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# 2. Fill in the new column's values.
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 4.5999999999999996
$ws.Range("Q6").Value = 4.2
$ws.Range("Q7").Value = 1.3
$ws.Range("Q8").Value = 10.8
$ws.Range("Q9").Value = 6.5
$ws.Range("Q10").Value = 2.9
$ws.Range("Q11").Value = 2.6
$ws.Range("Q12").Value = 13.1
$ws.Range("Q13").Value = 1
$ws.Range("Q14").Value = 1.3

# 3. Select column T (full-column selection), matching the recorded
#    sheetView selection in the target workbook.
$ws.Columns("T").Select()
